$wb = $excel.ActiveWorkbook

# OFF sheet - Week 17 (row 3, "R" = rest-of-season / sim totals) updated
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 208
$wsOff.Range("C3").Value = 137
$wsOff.Range("D3").Value = 40
$wsOff.Range("E3").Value = 15

# DEF sheet - Week 17 (row 3, "R") updated
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 197
$wsDef.Range("C3").Value = 142
$wsDef.Range("D3").Value = 48
$wsDef.Range("E3").Value = 25
